$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF") with the same style as the
# existing header cell H1 (bold, bordered, centered/top-aligned).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add data cells I2 and J2 with numeric values, matching the plain
# (unstyled) formatting of the other row-2 data cells.
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
